$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Nombre des joueurs (maximum 4 joueurs)" -> "Nombre des joueurs(2)"
#    Locate the paragraph (keeps the leading "N" run untouched, only
#    rewrites the remainder of the paragraph).
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Nombre des joueurs \(maximum 4 joueurs\)") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    # Replace everything after the leading "N" (keeps that run intact)
    # with the new wording, up to (but excluding) the paragraph mark.
    $body = $d.Range($pStart + 1, $pEnd - 1)
    $body.Text = "ombre des joueurs(2)"

    # ---------------------------------------------------------------
    # 2) Move the "_GoBack" bookmark so it now sits right after this
    #    paragraph's text (it previously sat after "Choix des pions").
    #    Bookmarks.Add on a truly collapsed (zero-length) range snaps
    #    to the start of the paragraph in this host, so we temporarily
    #    insert a 1-character placeholder at the desired spot, wrap the
    #    bookmark around it (non-zero width => exact placement), then
    #    delete the placeholder. The bookmark collapses in place and
    #    stays anchored exactly where we want it.
    # ---------------------------------------------------------------
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    $newEnd = $target.Range.End
    $insPoint = $d.Range($newEnd - 1, $newEnd - 1)
    $insPoint.InsertAfter("Z")

    $afterInsEnd = $target.Range.End
    $placeholder = $d.Range($afterInsEnd - 2, $afterInsEnd - 1)
    $d.Bookmarks.Add("_GoBack", $placeholder)

    $placeholder2 = $d.Range($afterInsEnd - 2, $afterInsEnd - 1)
    $placeholder2.Delete()
}
